# Update the "startup" sheet so that the SamplesTab and FilesTab rows
# reference the same Neo4jData/WebData file names as the CasesTab row
# (TC10_... instead of the stale TC01_... strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$ws.Range("D3").Value = "TC10_Canine_Filter_SamplePatho-TCellLymphoma_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC10_Canine_Filter_SamplePatho-TCellLymphoma_WebData.xlsx"

$ws.Range("D4").Value = "TC10_Canine_Filter_SamplePatho-TCellLymphoma_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC10_Canine_Filter_SamplePatho-TCellLymphoma_WebData.xlsx"

# Update the selected range to match the new selection in the saved file.
$ws.Range("D4:F4").Select()
